$d = $word.ActiveDocument

# Locate the Heading2 paragraph containing the "RUT" label
# (the book-code heading right before the italic "Ruth" sub-heading).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "RUT") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # The paragraph immediately following "RUT" holds only the italic
    # "Ruth" run (plus surrounding empty runs) and is to be removed
    # entirely, merging it out of the document.
    $toRemove = $target.Next()
    if ($toRemove -ne $null -and $toRemove.Range.Text.Trim() -eq "Ruth") {
        $toRemove.Range.Delete()
    }
}
